$d = $word.ActiveDocument

# Remove the existing _GoBack bookmark; it will be re-added in its new location
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$endRng = $d.Range($d.Content.End, $d.Content.End)

$rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>'
$rPrRed = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="FF0000"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>'

function PPr($ilvl) {
    return '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="' + $ilvl + '"/><w:numId w:val="2"/></w:numPr>' + $rPr + '</w:pPr>'
}

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' +
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
'<w:body>' +

'<w:p>' + (PPr 0) +
  '<w:r>' + $rPr + '<w:t>Boss Level Challenge 3 - Bitcoin Ticker</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
'</w:p>' +

'<w:p>' + (PPr 1) +
  '<w:r>' + $rPr + '<w:t>Widgets</w:t></w:r>' +
'</w:p>' +

'<w:p>' + (PPr 2) +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r>' + $rPr + '<w:t>DropdownButton</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r>' + $rPr + '<w:t>&lt;</w:t></w:r>' +
  '<w:r>' + $rPrRed + '<w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r>' + $rPr + '<w:t>&gt;</w:t></w:r>' +
'</w:p>' +

'<w:p>' + (PPr 2) +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r>' + $rPr + '<w:t>CupertinoPicker</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
'</w:p>' +

'<w:p>' + (PPr 1) +
  '<w:r>' + $rPr + '<w:t>Import</w:t></w:r>' +
'</w:p>' +

'<w:p>' + (PPr 2) +
  '<w:r>' + $rPr + '<w:t>Show Class ( just show this class in library).</w:t></w:r>' +
'</w:p>' +

'<w:p>' + (PPr 2) +
  '<w:r>' + $rPr + '<w:t>Hide Class (hide</w:t></w:r>' +
  '<w:r>' + $rPr + '<w:t xml:space="preserve"> this class in library).</w:t></w:r>' +
'</w:p>' +

'</w:body>' +
'</w:document>' +
'</pkg:xmlData></pkg:part>' +
'</pkg:package>'

$endRng.InsertXML($xml) | Out-Null
